# Fixed setting end dates for viranomaispaatos and kompostori
# - Adjusted tests: add a new "Uusi ilmoitus" (new notification) test row
#   for Karita Pyykoski to the ilmoitukset test fixture.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 2 (Johan Kemp's notification) into row 3 so the new row
# starts out with the same layout/formatting as the existing data row.
$ws.Rows("2:2").Copy()
$ws.Rows("3:3").PasteSpecial(-4104)

# Overwrite the copied values with the new notification's data
# (Karita Pyykoski).
$ws.Range("A3").Value2 = "20.6.2022"
$ws.Range("C3").Value2 = "Uusi ilmoitus"
$ws.Range("D3").Value2 = "Karita"
$ws.Range("E3").Value2 = "Pyykoski"
$ws.Range("F3").Value2 = "0400123645"
$ws.Range("G3").Value2 = "karita@pyykoski.fi"
$ws.Range("K3").Value2 = "Karita Pyykoski"
$ws.Range("O3").Value2 = "Kyykoski"
$ws.Range("AW3").Value2 = "Iivari"
$ws.Range("AX3").Value2 = "Kyykoski"

# Register the mail hyperlink for the new contact e-mail, same as G2.
$ws.Hyperlinks.Add($ws.Range("G3"), "mailto:karita@pyykoski.fi", [Type]::Missing, [Type]::Missing, "karita@pyykoski.fi")

# Match the saved selection in the edited workbook.
$ws.Range("O3").Select()
